$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1427.6666
$ws.Range("I34").Value = 1427.6666
$ws.Range("K34").Value = 1427.6666
$ws.Range("M34").Value = -1224.6666
$ws.Range("H36").Value = 1427.6666
$ws.Range("I36").Value = 1427.6666
$ws.Range("K36").Value = 1427.6666
$ws.Range("M36").Value = -712.6666
$ws.Range("H43").Value = 48149704
$ws.Range("J43").Value = 1997
$ws.Range("L43").Value = 1997
$ws.Range("N43").Value = -2135
$ws.Range("H51").Value = 71438570
$ws.Range("J51").Value = 250005000
$ws.Range("L51").Value = 250005000
$ws.Range("N51").Value = -250005968
$ws.Range("H86").Value = 78435520
$ws.Range("I86").Value = 76927224
$ws.Range("K86").Value = 76927224
$ws.Range("M86").Value = -76926101
$ws.Range("H89").Value = 78435520
$ws.Range("I89").Value = 76927224
$ws.Range("K89").Value = 384636120
$ws.Range("M89").Value = -384630504
$ws.Range("H109").Value = 56910.5
$ws.Range("I109").Value = 43821
$ws.Range("J109").Value = 70000
$ws.Range("K109").Value = 43821
$ws.Range("L109").Value = 70000
$ws.Range("M109").Value = -42434
$ws.Range("N109").Value = -72774
$ws.Range("H113").Value = 3999.3333
$ws.Range("I113").Value = 3999.3333
$ws.Range("K113").Value = 3999.3333
$ws.Range("M113").Value = -745.3332999999998
$ws.Range("H129").Value = 2500
$ws.Range("J129").Value = 2500
$ws.Range("L129").Value = 7500
$ws.Range("N129").Value = -17500
$ws.Range("H131").Value = 9171.625
$ws.Range("H132").Value = 1342.1777
$ws.Range("I132").Value = 1265.1395
$ws.Range("K132").Value = 3795.4185
$ws.Range("M132").Value = -1265.4185

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1504.5883
$ws.Range("I2").Value = 890.86957
$ws.Range("K2").Value = 890.86957
$ws.Range("M2").Value = -777.86957
$ws.Range("H45").Value = 56579.895
$ws.Range("I45").Value = 80433.234
$ws.Range("J45").Value = 4897.6665
$ws.Range("K45").Value = 80433.234
$ws.Range("L45").Value = 4897.6665
$ws.Range("M45").Value = -80056.234
$ws.Range("N45").Value = -5651.6665
$ws.Range("H61").Value = 967452.1
$ws.Range("I61").Value = 2969.4443
$ws.Range("K61").Value = 2969.4443
$ws.Range("M61").Value = -2757.4443
$ws.Range("H88").Value = 2243.2
$ws.Range("I88").Value = 2215
$ws.Range("K88").Value = 2215
$ws.Range("M88").Value = -1809
$ws.Range("H91").Value = 2243.2
$ws.Range("I91").Value = 2215
$ws.Range("K91").Value = 2215
$ws.Range("M91").Value = -811
$ws.Range("H92").Value = 77989.336
$ws.Range("J92").Value = 77989.336
$ws.Range("L92").Value = 77989.336
$ws.Range("N92").Value = -82981.336
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H116").Value = 1504.5883
$ws.Range("I116").Value = 890.86957
$ws.Range("K116").Value = 890.86957
$ws.Range("M116").Value = 1403.13043
$ws.Range("H132").Value = 4917
$ws.Range("I132").Value = 4829.6553
$ws.Range("K132").Value = 14488.9659
$ws.Range("M132").Value = -11958.9659
$ws.Range("H136").Value = 967452.1
$ws.Range("I136").Value = 2969.4443
$ws.Range("K136").Value = 8908.332900000001
$ws.Range("M136").Value = -6358.332900000001

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1504.5883
$ws.Range("I3").Value = 890.86957
$ws.Range("K3").Value = 890.86957
$ws.Range("M3").Value = -776.86957
$ws.Range("H137").Value = 103370.336
$ws.Range("J137").Value = 105555.5
$ws.Range("L137").Value = 105555.5
$ws.Range("N137").Value = -115755.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2660.9814
$ws.Range("I31").Value = 2241.8333
$ws.Range("K31").Value = 2241.8333
$ws.Range("M31").Value = -1946.8333
$ws.Range("H34").Value = 2660.9814
$ws.Range("I34").Value = 2241.8333
$ws.Range("K34").Value = 2241.8333
$ws.Range("M34").Value = -2039.8333
$ws.Range("H105").Value = 2437.077
$ws.Range("I105").Value = 1836.8889
$ws.Range("K105").Value = 1836.8889
$ws.Range("M105").Value = -89.88889999999992
$ws.Range("H107").Value = 2462.625
$ws.Range("I107").Value = 2304.4167
$ws.Range("K107").Value = 2304.4167
$ws.Range("M107").Value = -384.4167000000002
$ws.Range("H120").Value = 25000
$ws.Range("J120").Value = 25000
$ws.Range("L120").Value = 25000
$ws.Range("N120").Value = -32258

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1000372
$ws.Range("I4").Value = 1148462.2
$ws.Range("K4").Value = 3445386.6
$ws.Range("M4").Value = -3445274.6
$ws.Range("H11").Value = 685.7692
$ws.Range("I11").Value = 659.5833
$ws.Range("K11").Value = 1978.7499
$ws.Range("M11").Value = -1838.7499
$ws.Range("H12").Value = 576.3684
$ws.Range("I12").Value = 81
$ws.Range("K12").Value = 243
$ws.Range("M12").Value = -70
$ws.Range("H19").Value = 5149.8335
$ws.Range("J19").Value = 6666.3335
$ws.Range("L19").Value = 19999.0005
$ws.Range("N19").Value = -20347.0005
$ws.Range("H92").Value = 550
$ws.Range("J92").Value = 300
$ws.Range("L92").Value = 900
$ws.Range("N92").Value = -3396
$ws.Range("H109").Value = 16732.783
$ws.Range("I109").Value = 5363
$ws.Range("K109").Value = 16089
$ws.Range("M109").Value = -15049
$ws.Range("H113").Value = 555.7368
$ws.Range("I113").Value = 996.5
$ws.Range("J113").Value = 503.88235
$ws.Range("K113").Value = 2989.5
$ws.Range("L113").Value = 1511.64705
$ws.Range("M113").Value = -819.5
$ws.Range("N113").Value = -5851.64705
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H120").Value = 24500
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H133").Value = 2798.7273
$ws.Range("I133").Value = 2798.7273
$ws.Range("K133").Value = 8396.1819
$ws.Range("M133").Value = -3336.1819
$ws.Range("H138").Value = 2829.4707
$ws.Range("I138").Value = 2829.4707
$ws.Range("K138").Value = 8488.4121
$ws.Range("M138").Value = -3348.4121

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 41777890
$ws.Range("I80").Value = 136639.67
$ws.Range("J80").Value = 166701630
$ws.Range("K80").Value = 136639.67
$ws.Range("L80").Value = 166701630
$ws.Range("M80").Value = -135641.67
$ws.Range("N80").Value = -166703626
$ws.Range("H83").Value = 41777890
$ws.Range("I83").Value = 136639.67
$ws.Range("J83").Value = 166701630
$ws.Range("K83").Value = 683198.3500000001
$ws.Range("L83").Value = 833508150
$ws.Range("M83").Value = -678206.3500000001
$ws.Range("N83").Value = -833518134

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1520.4166
$ws.Range("I16").Value = 1344.5
$ws.Range("K16").Value = 1344.5
$ws.Range("M16").Value = -1174.5
$ws.Range("H40").Value = 1995.2368
$ws.Range("I40").Value = 1168.8667
$ws.Range("J40").Value = 5094.125
$ws.Range("K40").Value = 1168.8667
$ws.Range("L40").Value = 5094.125
$ws.Range("M40").Value = -1032.8667
$ws.Range("N40").Value = -5366.125
$ws.Range("H82").Value = 1242.75
$ws.Range("I82").Value = 990.6667
$ws.Range("K82").Value = 990.6667
$ws.Range("M82").Value = -629.6667
$ws.Range("H85").Value = 1242.75
$ws.Range("I85").Value = 990.6667
$ws.Range("K85").Value = 990.6667
$ws.Range("M85").Value = 257.3333
$ws.Range("H87").Value = 1000000000
$ws.Range("J87").Value = 1000000000
$ws.Range("L87").Value = 1000000000
$ws.Range("N87").Value = -1000002246
$ws.Range("H88").Value = 26249.75
$ws.Range("I88").Value = 27500
$ws.Range("J88").Value = 24999.5
$ws.Range("K88").Value = 27500
$ws.Range("L88").Value = 24999.5
$ws.Range("M88").Value = -27072
$ws.Range("N88").Value = -25855.5
$ws.Range("H90").Value = 1000000000
$ws.Range("J90").Value = 1000000000
$ws.Range("L90").Value = 3000000000
$ws.Range("N90").Value = -3000011232
$ws.Range("H91").Value = 26249.75
$ws.Range("I91").Value = 27500
$ws.Range("J91").Value = 24999.5
$ws.Range("K91").Value = 27500
$ws.Range("L91").Value = 24999.5
$ws.Range("M91").Value = -26018
$ws.Range("N91").Value = -27963.5
$ws.Range("H122").Value = 3307.804
$ws.Range("I122").Value = 2885.853
$ws.Range("J122").Value = 4151.706
$ws.Range("K122").Value = 8657.559000000001
$ws.Range("L122").Value = 12455.118
$ws.Range("M122").Value = -6207.559000000001
$ws.Range("N122").Value = -17355.118
$ws.Range("H136").Value = 2599.726
$ws.Range("I136").Value = 2516.7036
$ws.Range("K136").Value = 7550.110799999999
$ws.Range("M136").Value = -5000.110799999999

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10400.777
$ws.Range("J41").Value = 10831
$ws.Range("L41").Value = 10831
$ws.Range("N41").Value = -11611
$ws.Range("H132").Value = 2021.7705
$ws.Range("I132").Value = 1555.4222
$ws.Range("J132").Value = 3333.375
$ws.Range("K132").Value = 4666.2666
$ws.Range("L132").Value = 10000.125
$ws.Range("M132").Value = -2136.2666
$ws.Range("N132").Value = -15060.125

Write-Host "Applied all cell updates."